$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.833.89"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "2.083.65"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.29"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.93"
$ws.Range("E7").Value = "  +2.85%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0790"
$ws.Range("E10").Value = "  +1.93%  "
$ws.Range("E11").Value = "  +2.93%  "
$ws.Range("D12").Value = "2.390.28"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.77"
$ws.Range("E13").Value = "  +2.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.24"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.767"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").Value = "2.091.76"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "37.737.54"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.35"
$ws.Range("E20").Value = "  +2.48%  "
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.02"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -1.28%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.38"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("E27").Value = "  +6.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +2.13%  "
$ws.Range("E31").Value = "  +1.70%  "
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.72"
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0631"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("E36").Value = "  +2.84%  "
$ws.Range("E37").Value = "  +2.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0992"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.82"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.44"
$ws.Range("E43").Value = "  +6.38%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").Value = "1.465.94"
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.02"
$ws.Range("E47").Value = "  +5.46%  "
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.41"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").Value = "2.275.14"
$ws.Range("E51").Value = "  +0.42%  "
